$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.739.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.526.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.525.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.49"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.948.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.65"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.838.37"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.506.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.79%  "
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("E25").Value = "  -10.17%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0796"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.44"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.95"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "315.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.839"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0534"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0946"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.74"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.16%  "
